$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.518.38'
$ws.Range("E2").Value = '  +1.31%  '

$ws.Range("D3").Value = '1.878.37'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.025'
$ws.Range("E4").Value = '  +2.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.78'
$ws.Range("E5").Value = '  +1.40%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.023'
$ws.Range("E6").Value = '  +1.92%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5151'
$ws.Range("E7").Value = '  +0.55%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3949'
$ws.Range("E8").Value = '  +0.83%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08346'
$ws.Range("E9").Value = '  +0.27%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.119'
$ws.Range("E10").Value = '  +0.01%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.11'
$ws.Range("E11").Value = '  +1.57%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.258'
$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.45'
$ws.Range("E13").Value = '  -0.94%  '

$ws.Range("D14").Value = '1.848.98'
$ws.Range("E14").Value = '  -0.34%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.250'
$ws.Range("E15").Value = '  -0.16%  '

$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.024'
$ws.Range("E16").Value = '  +2.20%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001110'
$ws.Range("E17").Value = '  +0.92%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '91.54'
$ws.Range("E18").Value = '  +0.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06777'
$ws.Range("E19").Value = '  +2.03%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.023'
$ws.Range("E20").Value = '  +1.98%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.70'
$ws.Range("E21").Value = '  -0.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.976'
$ws.Range("E22").Value = '  -0.72%  '

$ws.Range("D23").Value = '28.560.88'
$ws.Range("E23").Value = '  +1.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.20'
$ws.Range("E24").Value = '  +0.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.267'
$ws.Range("E25").Value = '  +0.13%  '

$ws.Range("B26").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C26").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D26").Value = '2.075.99'
$ws.Range("E26").Value = '  -0.50%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '161.96'
$ws.Range("E27").Value = '  +2.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.72'
$ws.Range("E28").Value = '  +0.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.382'
$ws.Range("E29").Value = '  -4.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.43'
$ws.Range("E30").Value = '  +1.98%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1046'
$ws.Range("E31").Value = '  -1.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.036'
$ws.Range("E32").Value = '  -0.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.843'
$ws.Range("E33").Value = '  -0.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.657'
$ws.Range("E34").Value = '  +1.73%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02437'

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06502'
$ws.Range("E36").Value = '  -0.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '9.137'
$ws.Range("E37").Value = '  -6.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2186'
$ws.Range("E38").Value = '  +0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.253'
$ws.Range("E39").Value = '  +1.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.190'
$ws.Range("E40").Value = '  -1.12%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6450'
$ws.Range("E41").Value = '  -0.75%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.003'
$ws.Range("E42").Value = '  +0.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.24'
$ws.Range("E43").Value = '  -0.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6040'
$ws.Range("E44").Value = '  -1.38%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.03'
$ws.Range("E45").Value = '  -0.39%  '

$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.718'
$ws.Range("E46").Value = '  +1.27%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.255'
$ws.Range("E47").Value = '  -2.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.998'
$ws.Range("E48").Value = '  -0.47%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.214'
$ws.Range("E49").Value = '  -1.38%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.23'
$ws.Range("E50").Value = '  +0.93%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06871'
$ws.Range("E51").Value = '  -0.22%  '
